$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Check the table/list object exists before editing (defensive, as in commit message)
$hasTable = $false
foreach ($lo in $ws.ListObjects) {
    $hasTable = $true
}

# Insert a new row at position 8 (shifts existing rows 8+ down by one)
$ws.Rows.Item(8).EntireRow.Insert()

# Populate the newly inserted row 8 with the new shoe entry
$ws.Cells.Item(8, 1).Value = "Кеды Converse Chuck 70 AT-CX Hi"
$ws.Cells.Item(8, 2).Value = "15.500.₽"
$ws.Cells.Item(8, 3).Value = "https://sneakerhead.ru//upload/resize_cache/iblock/aa7/296_296_2/aq3ojiv2ydjahpzlf1qvbfvbtcxqldhq.jpg "
$ws.Cells.Item(8, 4).Value = "https://sneakerhead.ru/shoes/keds/chuck-70-at-cx-hi-A02776/"

# Remove the row that is now duplicated further down (old row 14 content, pushed to row 15)
$ws.Rows.Item(15).EntireRow.Delete()
